$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 3 corresponds to the 889e6598 file.
# Only this row's handoff/handback datetimes are regenerated.
$wsZhCn.Range("D3").Value = "2016-01-18 05:44:19"
$wsZhCn.Range("G3").Value = "2016-01-18 05:45:10"

# de-de sheet: row 3 corresponds to the 889e6598 file.
$wsDeDe.Range("D3").Value = "2016-01-18 05:44:32"
$wsDeDe.Range("G3").Value = "2016-01-18 05:45:32"
